$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.864.50'
$ws.Range('E2').Value = '  +0.35%  '

$ws.Range('D3').Value = '2.534.53'
$ws.Range('E3').Value = '  +0.55%  '

$ws.Range('E4').Value = '  +0.02%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '317.96'
$ws.Range('D5').Style = $ws.Range('B5').Style
$ws.Range('E5').Value = '  +4.54%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '95.69'
$ws.Range('D6').Style = $ws.Range('B6').Style
$ws.Range('E6').Value = '  -1.57%  '

$ws.Range('E7').Value = '  +0.66%  '

$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.533'
$ws.Range('D9').Style = $ws.Range('B9').Style
$ws.Range('E9').Value = '  -1.25%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.31'
$ws.Range('D10').Style = $ws.Range('B10').Style
$ws.Range('E10').Value = '  -1.16%  '

$ws.Range('E11').Value = '  -0.05%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.63'
$ws.Range('D12').Style = $ws.Range('B12').Style
$ws.Range('E12').Value = '  -0.56%  '

$ws.Range('E13').Value = '  -0.57%  '

$ws.Range('D14').Value = '2.926.73'
$ws.Range('E14').Value = '  +0.70%  '

$ws.Range('D15').Value = '2.551.21'
$ws.Range('E15').Value = '  +1.42%  '

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '15.36'
$ws.Range('D16').Style = $ws.Range('B16').Style
$ws.Range('E16').Value = '  +2.20%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.850'
$ws.Range('D17').Style = $ws.Range('B17').Style
$ws.Range('E17').Value = '  -0.46%  '

$ws.Range('D18').Value = '42.884.86'
$ws.Range('E18').Value = '  +0.40%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.05'
$ws.Range('D19').Style = $ws.Range('B19').Style
$ws.Range('E19').Value = '  +0.94%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.64'
$ws.Range('D20').Style = $ws.Range('B20').Style
$ws.Range('E20').Value = '  +3.24%  '

$ws.Range('D21').Value = '0.0₃0970'
$ws.Range('E21').Value = '  -0.34%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '70.21'
$ws.Range('D22').Style = $ws.Range('B22').Style
$ws.Range('E22').Value = '  -1.33%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '252.14'
$ws.Range('D23').Style = $ws.Range('B23').Style
$ws.Range('E23').Value = '  +0.56%  '

$ws.Range('E24').Value = '  +1.79%  '

$ws.Range('E25').Value = '  -0.01%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '27.05'
$ws.Range('D26').Style = $ws.Range('B26').Style
$ws.Range('E26').Value = '  +0.62%  '

$ws.Range('E27').Value = '  +0.02%  '

$ws.Range('E28').Value = '  +4.21%  '

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '40.00'
$ws.Range('D29').Style = $ws.Range('B29').Style
$ws.Range('E29').Value = '  +5.41%  '

$ws.Range('E30').Value = '  -0.59%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.09'
$ws.Range('D31').Style = $ws.Range('B31').Style
$ws.Range('E31').Value = '  +1.77%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '153.75'
$ws.Range('D32').Style = $ws.Range('B32').Style
$ws.Range('E32').Value = '  -1.68%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.12'
$ws.Range('D33').Style = $ws.Range('B33').Style
$ws.Range('E33').Value = '  +2.18%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.32'
$ws.Range('D34').Style = $ws.Range('B34').Style

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '19.02'
$ws.Range('D35').Style = $ws.Range('B35').Style
$ws.Range('E35').Value = '  +4.07%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0791'
$ws.Range('D36').Style = $ws.Range('B36').Style
$ws.Range('E36').Value = '  +0.13%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.62'
$ws.Range('D37').Style = $ws.Range('B37').Style
$ws.Range('E37').Value = '  +0.17%  '

$ws.Range('E38').Value = '  -3.20%  '

$ws.Range('E39').Value = '  +0.00%  '

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '23.96'
$ws.Range('D40').Style = $ws.Range('B40').Style
$ws.Range('E40').Value = '  -0.20%  '

$ws.Range('E41').Value = '  +10.22%  '

$ws.Range('E42').Value = '  -0.70%  '

$ws.Range('E43').Value = '  +1.51%  '

$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = $ws.Range('B44').Style
$ws.Range('E44').Value = '  +0.37%  '

$ws.Range('B45').Value = 'NEARProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.33'
$ws.Range('D45').Style = $ws.Range('B45').Style
$ws.Range('E45').Value = '  -1.64%  '

$ws.Range('D46').Value = '2.021.69'
$ws.Range('E46').Value = '  -0.34%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '85.57'
$ws.Range('D47').Style = $ws.Range('B47').Style
$ws.Range('E47').Value = '  +0.34%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.80'
$ws.Range('D48').Style = $ws.Range('B48').Style
$ws.Range('E48').Value = '  -2.02%  '

$ws.Range('D49').Value = '2.784.73'
$ws.Range('E49').Value = '  +0.69%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '74.46'
$ws.Range('D50').Style = $ws.Range('B50').Style
$ws.Range('E50').Value = '  +3.10%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '102.41'
$ws.Range('D51').Style = $ws.Range('B51').Style
$ws.Range('E51').Value = '  +0.67%  '
